# "Generate Report for Archive"
#
# The localization status changes from "Ready for handoff" to
# "In Translation" for both tracked files, on every sheet that surfaces
# that status:
#   - Overview : columns "zh-cn" (E) and "de-de" (F), rows 2-3
#   - zh-cn    : "Status" column (C), rows 2-3
#   - de-de    : "Status" column (C), rows 2-3
#
# Because the new text is shorter than the old text, the status columns
# are narrowed to fit on Overview/zh-cn/de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status cells, rows 2-3 ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), rows 2-3 ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C), rows 2-3 ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Narrow the now-shorter status columns to fit the new text ---
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
